$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "MIE(9.228565317434532, 5.102986966614236, -9.369948010795211, 11.196645443370722)"
$ws.Range("C2").Value = "NIG(1.5545601925557169, 1.253905792980563, 3.5441046053262584, 5.5370264167077385)"
$ws.Range("D2").Value = "NCT(2.784755684251403, 1.484074419609915, -0.8403187444912417, 2.642922373263103)"
$ws.Range("E2").Value = "NIG(1.0734894235987824, 0.7950558664099221, 4.611890505775763, 6.752311911919691)"
